$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 328.18182
$ws.Range("I33").Value = 341.5
$ws.Range("J33").Value = 195
$ws.Range("K33").Value = 341.5
$ws.Range("L33").Value = 195
$ws.Range("M33").Value = -112.5
$ws.Range("N33").Value = -653
$ws.Range("H41").Value = 2803.4
$ws.Range("I41").Value = 2500
$ws.Range("K41").Value = 2500
$ws.Range("M41").Value = -2060
$ws.Range("H53").Value = 464.54544
$ws.Range("I53").Value = 473.75
$ws.Range("K53").Value = 473.75
$ws.Range("M53").Value = 163.25
$ws.Range("H76").Value = 4999.3335
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 4999.3335
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 4999.3335
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -5629.3335
$ws.Range("H79").Value = 4999.3335
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 4999.3335
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 4999.3335
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -7183.3335
$ws.Range("H113").Value = 27417.092
$ws.Range("I113").Value = 27707.375
$ws.Range("J113").Value = 26643
$ws.Range("K113").Value = 27707.375
$ws.Range("L113").Value = 26643
$ws.Range("M113").Value = -24453.375
$ws.Range("N113").Value = -33151
$ws.Range("H141").Value = 7238.6665
$ws.Range("I141").Value = 7137.8667
$ws.Range("K141").Value = 21413.6001
$ws.Range("M141").Value = -16233.6001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19189.562
$ws.Range("I32").Value = 17387.322
$ws.Range("K32").Value = 17387.322
$ws.Range("M32").Value = -17100.322
$ws.Range("H46").Value = 6050
$ws.Range("I46").Value = 5985.2
$ws.Range("J46").Value = 6114.8
$ws.Range("K46").Value = 5985.2
$ws.Range("L46").Value = 6114.8
$ws.Range("M46").Value = -5666.2
$ws.Range("N46").Value = -6752.8
$ws.Range("H132").Value = 7811.423
$ws.Range("I132").Value = 5168.684
$ws.Range("K132").Value = 15506.052
$ws.Range("M132").Value = -12976.052
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 21966.666
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 21966.666
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 21966.666
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -22596.666
$ws.Range("H79").Value = 21966.666
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 21966.666
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 21966.666
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -24150.666
$ws.Range("H105").Value = 2068.7896
$ws.Range("I105").Value = 2220.7334
$ws.Range("K105").Value = 2220.7334
$ws.Range("M105").Value = -473.7334000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 257870
$ws.Range("I31").Value = 401631.25
$ws.Range("K31").Value = 401631.25
$ws.Range("M31").Value = -401336.25
$ws.Range("H34").Value = 257870
$ws.Range("I34").Value = 401631.25
$ws.Range("K34").Value = 401631.25
$ws.Range("M34").Value = -401429.25
$ws.Range("H58").Value = 2528.125
$ws.Range("I58").Value = 2733.35
$ws.Range("J58").Value = 1502
$ws.Range("K58").Value = 2733.35
$ws.Range("L58").Value = 1502
$ws.Range("M58").Value = -2530.35
$ws.Range("N58").Value = -1908
$ws.Range("H62").Value = 4391.0713
$ws.Range("I62").Value = 3996.8
$ws.Range("K62").Value = 3996.8
$ws.Range("M62").Value = -3372.8
$ws.Range("H65").Value = 4391.0713
$ws.Range("I65").Value = 3996.8
$ws.Range("K65").Value = 19984
$ws.Range("M65").Value = -16864
$ws.Range("H88").Value = 11624
$ws.Range("J88").Value = 13498.5
$ws.Range("L88").Value = 13498.5
$ws.Range("N88").Value = -14310.5
$ws.Range("H91").Value = 11624
$ws.Range("J91").Value = 13498.5
$ws.Range("L91").Value = 13498.5
$ws.Range("N91").Value = -16306.5
$ws.Range("H99").Value = 1113817.1
$ws.Range("I99").Value = 1669243.1
$ws.Range("K99").Value = 1669243.1
$ws.Range("M99").Value = -1667745.1
$ws.Range("H122").Value = 1443.52
$ws.Range("I122").Value = 1106.1428
$ws.Range("K122").Value = 3318.4284
$ws.Range("M122").Value = -868.4284000000002
$ws.Range("H126").Value = 1113817.1
$ws.Range("I126").Value = 1669243.1
$ws.Range("K126").Value = 5007729.300000001
$ws.Range("M126").Value = -5005259.300000001
$ws.Range("H132").Value = 3897.0278
$ws.Range("I132").Value = 2579.8
$ws.Range("J132").Value = 50000
$ws.Range("K132").Value = 7739.400000000001
$ws.Range("L132").Value = 150000
$ws.Range("M132").Value = -5209.400000000001
$ws.Range("N132").Value = -155060
$ws.Range("H136").Value = 2528.125
$ws.Range("I136").Value = 2733.35
$ws.Range("J136").Value = 1502
$ws.Range("K136").Value = 8200.049999999999
$ws.Range("L136").Value = 4506
$ws.Range("M136").Value = -5650.049999999999
$ws.Range("N136").Value = -9606
$ws.Range("H141").Value = 119218.805
$ws.Range("J141").Value = 119218.805
$ws.Range("L141").Value = 119218.805
$ws.Range("N141").Value = -129578.805
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1856.7142
$ws.Range("I80").Value = 1832.8334
$ws.Range("K80").Value = 5498.5002
$ws.Range("M80").Value = -4562.5002
$ws.Range("H83").Value = 1856.7142
$ws.Range("I83").Value = 1832.8334
$ws.Range("K83").Value = 16495.5006
$ws.Range("M83").Value = -11815.5006
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 63933.21
$ws.Range("I132").Value = 65390.188
$ws.Range("J132").Value = 56162.668
$ws.Range("K132").Value = 196170.564
$ws.Range("L132").Value = 168488.004
$ws.Range("M132").Value = -193640.564
$ws.Range("N132").Value = -173548.004
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8444.214
$ws.Range("I40").Value = 9246.166999999999
$ws.Range("J40").Value = 3632.5
$ws.Range("K40").Value = 9246.166999999999
$ws.Range("L40").Value = 3632.5
$ws.Range("M40").Value = -9110.166999999999
$ws.Range("N40").Value = -3904.5
$ws.Range("H55").Value = 10830.333
$ws.Range("I55").Value = 804.5833
$ws.Range("J55").Value = 50933.332
$ws.Range("K55").Value = 804.5833
$ws.Range("L55").Value = 50933.332
$ws.Range("M55").Value = -631.5833
$ws.Range("N55").Value = -51279.332
$ws.Range("H100").Value = 4923.2856
$ws.Range("I100").Value = 4746.5
$ws.Range("K100").Value = 4746.5
$ws.Range("M100").Value = -4205.5
$ws.Range("H132").Value = 7518.8477
$ws.Range("I132").Value = 3880.9412
$ws.Range("J132").Value = 9651.414000000001
$ws.Range("K132").Value = 11642.8236
$ws.Range("L132").Value = 28954.242
$ws.Range("M132").Value = -9112.8236
$ws.Range("N132").Value = -34014.242
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 67972.664
$ws.Range("J46").Value = 67972.664
$ws.Range("L46").Value = 67972.664
$ws.Range("N46").Value = -68434.664
$ws.Range("H81").Value = 3153.36
$ws.Range("I81").Value = 3551.4736
$ws.Range("K81").Value = 7102.9472
$ws.Range("M81").Value = -6041.9472
$ws.Range("H84").Value = 3153.36
$ws.Range("I84").Value = 3551.4736
$ws.Range("K84").Value = 35514.736
$ws.Range("M84").Value = -30210.736
$ws.Range("H96").Value = 95724.17999999999
$ws.Range("I96").Value = 203471.4
$ws.Range("K96").Value = 203471.4
$ws.Range("M96").Value = -202098.4
$ws.Range("H126").Value = 1657.4
$ws.Range("I126").Value = 1134
$ws.Range("J126").Value = 2629.4285
$ws.Range("K126").Value = 3402
$ws.Range("L126").Value = 7888.2855
$ws.Range("M126").Value = -932
$ws.Range("N126").Value = -12828.2855
$ws.Range("H130").Value = 37476.332
$ws.Range("J130").Value = 37476.332
$ws.Range("L130").Value = 37476.332
$ws.Range("N130").Value = -47516.332
$ws.Range("H132").Value = 2450.7742
$ws.Range("I132").Value = 2413.5
$ws.Range("J132").Value = 2578.5715
$ws.Range("K132").Value = 7240.5
$ws.Range("L132").Value = 7735.7145
$ws.Range("M132").Value = -4710.5
$ws.Range("N132").Value = -12795.7145
$ws.Range("H134").Value = 67972.664
$ws.Range("J134").Value = 67972.664
$ws.Range("L134").Value = 203917.992
$ws.Range("N134").Value = -208987.992
